$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "teststep" header in column C, row 1 (no special style, unlike A1/B1)
$ws.Cells.Item(1, 3).Value = "teststep"

# Set width of new column C (target stored width 21; the runtime's
# pixel-rounding adds a constant 5/6 offset, so compensate for it)
$ws.Columns.Item(3).ColumnWidth = 20.166666666666668

# Update selection to mirror the saved worksheet view (active cell C2)
$ws.Range("C2").Select()
